# Update "想去人数" (F column) figures across the 展览, 演出 and 全部类型 sheets
# to reflect the latest scrape output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    2  = 5936
    4  = 1177
    5  = 1091
    6  = 863
    10 = 71
    11 = 40
    12 = 31
    13 = 2144
    14 = 1540
    15 = 1197
    17 = 218
    18 = 472
    19 = 696
    20 = 254
    23 = 526
    24 = 3979
    25 = 203
    26 = 138
    28 = 181
    29 = 65
    30 = 578
    31 = 17
    35 = 343
    36 = 894
    38 = 82
    39 = 101
    40 = 98
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$sheet2Updates = @{
    3 = 767
}
foreach ($row in $sheet2Updates.Keys) {
    $ws2.Cells.Item($row, 6).Value = $sheet2Updates[$row]
}

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 5936
    4  = 1177
    6  = 767
    7  = 1091
    8  = 863
    14 = 71
    15 = 40
    17 = 31
    18 = 2144
    19 = 1540
    20 = 1197
    22 = 218
    23 = 472
    25 = 696
    26 = 254
    29 = 526
    30 = 3979
    31 = 203
    32 = 138
    34 = 181
    35 = 65
    36 = 578
    37 = 17
    41 = 343
    42 = 894
    44 = 82
    45 = 101
    46 = 98
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
